# Generate Report for Handback
#
# A new handback run produced a fresh pair of GUID-named files (and a new
# content hash / timestamps) for both the zh-cn and de-de locales. Update
# the Overview / zh-cn / de-de sheets so every cell and hyperlink that
# referenced the old handoff/handback file names now reflects the new
# ones, and the correspond handoff/handback datetimes are refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New values produced by this handback run
# ---------------------------------------------------------------------
$oldMd1 = "c12cd017-6b27-4c1a-a737-cf9612fd3db3.md"
$newMd1 = "732f54ce-d247-41d5-9667-659db39b894b.md"

$oldMd2 = "dd08a102-2d58-4b5b-bcf7-e8440b9b0360.md"
$newMd2 = "ffff8c528445-32b9-46dc-8af7-c8db5621ef18.md"

$newXlfZh = "732f54ce-d247-41d5-9667-659db39b894b.2d58cbefe895a089d37b989a395bcb42b22e00a8.zh-cn.xlf"
$newXlfDe = "732f54ce-d247-41d5-9667-659db39b894b.2d58cbefe895a089d37b989a395bcb42b22e00a8.de-de.xlf"

$newHandoffZh = "2016-03-13 07:04:31"
$newHandbackZh = "2016-03-13 07:04:49"
$newHandoffDe = "2016-03-13 07:04:35"
$newHandbackDe = "2016-03-13 07:04:54"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("A3").Value = $newMd2

$wsOverview.Range("A2").Hyperlinks.Delete()

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/530a43617073255f24d34aa9a109b70089ae3c44/e2e/$oldMd1",
    "",
    "",
    $newMd1
) | Out-Null
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/530a43617073255f24d34aa9a109b70089ae3c44/e2e/$oldMd2",
    "",
    "",
    $newMd2
) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMd1
$wsZh.Range("D2").Value = $newXlfZh
$wsZh.Range("E2").Value = $newHandoffZh
$wsZh.Range("F2").Value = $newMd1
$wsZh.Range("G2").Value = $newXlfZh
$wsZh.Range("H2").Value = $newHandbackZh

$wsZh.Range("A3").Value = $newMd2
$wsZh.Range("D3").Value = $newXlfZh
$wsZh.Range("E3").Value = $newHandoffZh
$wsZh.Range("F3").Value = $newMd2
$wsZh.Range("G3").Value = $newXlfZh
$wsZh.Range("H3").Value = $newHandbackZh

$wsZh.Range("A2").Hyperlinks.Delete()

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/530a43617073255f24d34aa9a109b70089ae3c44/e2e/$oldMd1", "", "", $newMd1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/530a43617073255f24d34aa9a109b70089ae3c44/e2e/$oldMd1", "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5893dae372a1002971f4504cd993d9cf01678edd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c12cd017-6b27-4c1a-a737-cf9612fd3db3.d4dcba8639963b2215a27b38f8a847e51f789549.zh-cn.xlf", "", "", $newXlfZh) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ddbccd1f85423a9dc5f1eb4e22760e9801b2f9c4/e2e/$oldMd1", "", "", $newMd1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e4fd7bc9fed8720dbc9a8a65d9d526ef720b8ab8/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c12cd017-6b27-4c1a-a737-cf9612fd3db3.d4dcba8639963b2215a27b38f8a847e51f789549.zh-cn.xlf", "", "", $newXlfZh) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/530a43617073255f24d34aa9a109b70089ae3c44/e2e/$oldMd2", "", "", $newMd2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/530a43617073255f24d34aa9a109b70089ae3c44/e2e/$oldMd2", "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5893dae372a1002971f4504cd993d9cf01678edd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/dd08a102-2d58-4b5b-bcf7-e8440b9b0360.ce40297b13a84debd09da7a3288f9090003c73ad.zh-cn.xlf", "", "", $newXlfZh) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ddbccd1f85423a9dc5f1eb4e22760e9801b2f9c4/e2e/$oldMd2", "", "", $newMd2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e4fd7bc9fed8720dbc9a8a65d9d526ef720b8ab8/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/dd08a102-2d58-4b5b-bcf7-e8440b9b0360.ce40297b13a84debd09da7a3288f9090003c73ad.zh-cn.xlf", "", "", $newXlfZh) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMd1
$wsDe.Range("D2").Value = $newXlfDe
$wsDe.Range("E2").Value = $newHandoffDe
$wsDe.Range("F2").Value = $newMd1
$wsDe.Range("G2").Value = $newXlfDe
$wsDe.Range("H2").Value = $newHandbackDe

$wsDe.Range("A3").Value = $newMd2
$wsDe.Range("D3").Value = $newXlfDe
$wsDe.Range("E3").Value = $newHandoffDe
$wsDe.Range("F3").Value = $newMd2
$wsDe.Range("G3").Value = $newXlfDe
$wsDe.Range("H3").Value = $newHandbackDe

$wsDe.Range("A2").Hyperlinks.Delete()

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/530a43617073255f24d34aa9a109b70089ae3c44/e2e/$oldMd1", "", "", $newMd1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/530a43617073255f24d34aa9a109b70089ae3c44/e2e/$oldMd1", "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d0f13483b83197341f9d73df040d3d811b4f28b6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c12cd017-6b27-4c1a-a737-cf9612fd3db3.d4dcba8639963b2215a27b38f8a847e51f789549.de-de.xlf", "", "", $newXlfDe) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/bba7b835cd161483ff833035eecfcb50151fd1fc/e2e/$oldMd1", "", "", $newMd1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f69da95fd7089b38a7e9909d24e4b00157c7882b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c12cd017-6b27-4c1a-a737-cf9612fd3db3.d4dcba8639963b2215a27b38f8a847e51f789549.de-de.xlf", "", "", $newXlfDe) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/530a43617073255f24d34aa9a109b70089ae3c44/e2e/$oldMd2", "", "", $newMd2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/530a43617073255f24d34aa9a109b70089ae3c44/e2e/$oldMd2", "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d0f13483b83197341f9d73df040d3d811b4f28b6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/dd08a102-2d58-4b5b-bcf7-e8440b9b0360.ce40297b13a84debd09da7a3288f9090003c73ad.de-de.xlf", "", "", $newXlfDe) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/bba7b835cd161483ff833035eecfcb50151fd1fc/e2e/$oldMd2", "", "", $newMd2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f69da95fd7089b38a7e9909d24e4b00157c7882b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/dd08a102-2d58-4b5b-bcf7-e8440b9b0360.ce40297b13a84debd09da7a3288f9090003c73ad.de-de.xlf", "", "", $newXlfDe) | Out-Null

"Generated handback report"
